# Auto-generated edit script applying scheduled market-data refresh values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 274.57144
$ws.Range("I28").Value = 322.18182
$ws.Range("K28").Value = 322.18182
$ws.Range("M28").Value = 162.81818
$ws.Range("H33").Value = 9657395
$ws.Range("J33").Value = 4500.1665
$ws.Range("L33").Value = 4500.1665
$ws.Range("N33").Value = -4958.1665
$ws.Range("H112").Value = 32987.85
$ws.Range("I112").Value = 2214.6667
$ws.Range("J112").Value = 39582.105
$ws.Range("K112").Value = 6644.000100000001
$ws.Range("L112").Value = 118746.315
$ws.Range("M112").Value = -5536.000100000001
$ws.Range("N112").Value = -120962.315
$ws.Range("H129").Value = 2497.8572
$ws.Range("H137").Value = 16667700
$ws.Range("J137").Value = 33334400
$ws.Range("L137").Value = 100003200
$ws.Range("N137").Value = -100008300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 223992.16
$ws.Range("I32").Value = 257492.11
$ws.Range("J32").Value = 6242.5
$ws.Range("K32").Value = 257492.11
$ws.Range("L32").Value = 6242.5
$ws.Range("M32").Value = -257205.11
$ws.Range("N32").Value = -6816.5
$ws.Range("H61").Value = 819513.6
$ws.Range("I61").Value = 2235.6853
$ws.Range("K61").Value = 2235.6853
$ws.Range("M61").Value = -2023.6853
$ws.Range("H136").Value = 819513.6
$ws.Range("I136").Value = 2235.6853
$ws.Range("K136").Value = 6707.0559
$ws.Range("M136").Value = -4157.0559

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11336.333
$ws.Range("I105").Value = 12448.556
$ws.Range("K105").Value = 12448.556
$ws.Range("M105").Value = -10701.556
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 11908289
$ws.Range("I16").Value = 12989952
$ws.Range("K16").Value = 12989952
$ws.Range("M16").Value = -12989665
$ws.Range("H33").Value = 4645.8
$ws.Range("I33").Value = 1140.1428
$ws.Range("J33").Value = 12825.667
$ws.Range("K33").Value = 1140.1428
$ws.Range("L33").Value = 12825.667
$ws.Range("M33").Value = -761.1428000000001
$ws.Range("N33").Value = -13583.667
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H58").Value = 2133.0715
$ws.Range("I58").Value = 2590.1428
$ws.Range("J58").Value = 1676
$ws.Range("K58").Value = 2590.1428
$ws.Range("L58").Value = 1676
$ws.Range("M58").Value = -2387.1428
$ws.Range("N58").Value = -2082
$ws.Range("H93").Value = 6459.5
$ws.Range("I93").Value = 6459.5
$ws.Range("K93").Value = 6459.5
$ws.Range("M93").Value = -4587.5
$ws.Range("H107").Value = 1942.6316
$ws.Range("I107").Value = 1934
$ws.Range("K107").Value = 1934
$ws.Range("M107").Value = -14
$ws.Range("H113").Value = 11908289
$ws.Range("I113").Value = 12989952
$ws.Range("K113").Value = 12989952
$ws.Range("M113").Value = -12987782
$ws.Range("H132").Value = 18881.846
$ws.Range("I132").Value = 21005.256
$ws.Range("K132").Value = 63015.768
$ws.Range("M132").Value = -60485.768
$ws.Range("H136").Value = 2133.0715
$ws.Range("I136").Value = 2590.1428
$ws.Range("J136").Value = 1676
$ws.Range("K136").Value = 7770.428400000001
$ws.Range("L136").Value = 5028
$ws.Range("M136").Value = -5220.428400000001
$ws.Range("N136").Value = -10128

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1000
$ws.Range("J92").Value = 1200
$ws.Range("L92").Value = 3600
$ws.Range("N92").Value = -6096
$ws.Range("H113").Value = 463.9643
$ws.Range("J113").Value = 424.79166
$ws.Range("L113").Value = 1274.37498
$ws.Range("N113").Value = -5614.374980000001
$ws.Range("H122").Value = 4167834.5
$ws.Range("I122").Value = 8334147
$ws.Range("J122").Value = 1521.75
$ws.Range("K122").Value = 75007323
$ws.Range("L122").Value = 13695.75
$ws.Range("M122").Value = -75004873
$ws.Range("N122").Value = -18595.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2686.2354
$ws.Range("I122").Value = 2681.88
$ws.Range("K122").Value = 8045.64
$ws.Range("M122").Value = -5595.64
$ws.Range("H123").Value = 53999.5
$ws.Range("J123").Value = 53999.5
$ws.Range("L123").Value = 53999.5
$ws.Range("N123").Value = -58899.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 10000
$ws.Range("I45").Value = 10000
$ws.Range("K45").Value = 10000
$ws.Range("M45").Value = -9593
$ws.Range("H48").Value = 50000
$ws.Range("I48").Value = 50000
$ws.Range("K48").Value = 50000
$ws.Range("M48").Value = -49339
$ws.Range("H55").Value = 997.7143
$ws.Range("I55").Value = 771.1177
$ws.Range("J55").Value = 1347.909
$ws.Range("K55").Value = 771.1177
$ws.Range("L55").Value = 1347.909
$ws.Range("M55").Value = -598.1177
$ws.Range("N55").Value = -1693.909
$ws.Range("H61").Value = 1989.6774
$ws.Range("I61").Value = 1922.7
$ws.Range("J61").Value = 3999
$ws.Range("K61").Value = 1922.7
$ws.Range("L61").Value = 3999
$ws.Range("M61").Value = -1720.7
$ws.Range("N61").Value = -4403
$ws.Range("H113").Value = 1989.6774
$ws.Range("I113").Value = 1922.7
$ws.Range("J113").Value = 3999
$ws.Range("K113").Value = 1922.7
$ws.Range("L113").Value = 3999
$ws.Range("M113").Value = 247.3
$ws.Range("N113").Value = -8339
$ws.Range("H132").Value = 3181.0667
$ws.Range("I132").Value = 2803.3157
$ws.Range("K132").Value = 8409.947100000001
$ws.Range("M132").Value = -5879.947100000001
$ws.Range("H136").Value = 1935.3368
$ws.Range("I136").Value = 2326
$ws.Range("K136").Value = 6978
$ws.Range("M136").Value = -4428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 258698.5
$ws.Range("J3").Value = 258698.5
$ws.Range("L3").Value = 258698.5
$ws.Range("N3").Value = -258926.5
$ws.Range("H132").Value = 3110.8667
$ws.Range("I132").Value = 2970.0667
$ws.Range("K132").Value = 8910.2001
$ws.Range("M132").Value = -6380.2001
$ws.Range("H136").Value = 21827.225
$ws.Range("I136").Value = 26837.422
$ws.Range("J136").Value = 4519.273
$ws.Range("K136").Value = 80512.266
$ws.Range("L136").Value = 13557.819
$ws.Range("M136").Value = -77962.266
$ws.Range("N136").Value = -18657.819

Write-Host "Applied 163 cell updates across 8 sheets"